$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B2 (shared string "aaron") to "a"
$ws.Range("B2").Value = "a"

# Populate the new cell B3 with "b" (as entered from a form)
$ws.Range("B3").Value = "b"
